# Update execution error diagram
#
# On slide 1, the table named "Table 24" lists persons (one per row under a
# "persons : UniquePersonList" header). Its second row currently reads
# "Amy" and is being renamed to "Bob".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$targetShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $shape = $s.Shapes.Item($i)
    if ($shape.HasTable -and $shape.Name -eq "Table 24") {
        $targetShape = $shape
    }
}

$tbl = $targetShape.Table
$cell = $tbl.Rows.Item(2).Cells.Item(1)
$cell.Shape.TextFrame.TextRange.Text = "Bob"
